$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "State" column into hotel_info, between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheet tabs: review_info first, hotel_info second ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))
